# final edits to proposal before submitting
#
# 1) Split the "Ken Jennings's tenure" bullet into a "notable contestant's
#    tenures" bullet (same run formatting throughout, just re-worded).
# 2) Wrap the "jarchive" run with spell-check proofErr markers.
#
# Both edits are applied by replacing the whole containing paragraph with a
# reconstructed OOXML paragraph via Range.InsertXML - this lets us emit
# the exact run layout (including the <w:proofErr/> markers that Word's
# live spell-checker would normally add) instead of relying on inline
# Find/Replace, which cannot introduce non-text markup.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# Edit 1: "Complete a case study of Ken Jennings's tenure on Jeopardy"
#      -> "Complete a case study of notable contestant's tenures on Jeopardy"
#         split across five runs with identical formatting
# ---------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

$para1 = '<w:p w14:paraId="4DB645A2" w14:textId="3A04E966" w:rsidR="00F55A8D" w:rsidRDefault="00C6143F" w:rsidP="00F55A8D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">Complete a case study of </w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>notable contestant' + [char]0x2019 + 's</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> tenure</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>s</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> on Jeopardy</w:t></w:r>' + `
  '</w:p>'

$r1 = $d.Content
$r1.Find.Execute("Complete a case study of Ken Jennings" + [char]0x2019 + "s tenure on Jeopardy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target1 = $r1.Paragraphs(1).Range
$target1.InsertXML($pkgHeader + $para1 + $pkgFooter)

# ---------------------------------------------------------------------
# Edit 2: wrap the "jarchive" run in spellStart/spellEnd proofErr markers
#
# InsertXML only splices cleanly in place when the target range's end
# coincides with a paragraph boundary (otherwise it appends a sibling
# paragraph instead of substituting in place). So rather than rebuild the
# whole paragraph - which would require reconstructing the preceding
# hyperlink run, and the engine's OOXML importer drops unresolved
# <w:rStyle> references on any run that passes through InsertXML - only
# the tail starting just before "jarchive" (through the end of the
# paragraph) is rebuilt. Everything before that, including the hyperlink,
# is left completely untouched. (The single-space run immediately before
# "jarchive" is pulled into the rebuilt tail too, since splicing exactly
# at its trailing edge gets silently absorbed into the previous run by
# the importer.)
# ---------------------------------------------------------------------
$tail2 = '<w:r w:rsidR="005F30B8"><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r w:rsidRPr="005F30B8"><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t>jarchive</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r w:rsidRPr="005F30B8"><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="14"/><w:szCs w:val="14"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>script/package in R</w:t></w:r>' + `
  '<w:r w:rsidR="007811E5"><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> to collect additional data fields required for above objectives</w:t></w:r>'
$para2 = '<w:p>' + $tail2 + '</w:p>'

$r2 = $d.Content
$r2.Find.Execute("jarchive", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraEnd2 = $r2.Paragraphs(1).Range.End
$target2 = $d.Range($r2.Start - 1, $paraEnd2)
$target2.InsertXML($pkgHeader + $para2 + $pkgFooter)

Write-Host "Applied final edits."
